$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Semestre ideal:" value from "EM-5" to "EF-5,EM-5" (columns B and C)
$ws.Range("B9").Value = "EF-5,EM-5"
$ws.Range("C9").Value = "EF-5,EM-5"

# Remove the "Requisitos:" row (22) and the requirement text row (23)
$ws.Range("A22:C23").EntireRow.Delete()
